$d = $word.ActiveDocument

# Locate the paragraph that ends with "And use babel for compilation of component to js"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a`n") -eq "And use babel for compilation of component to js") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$following = $target.Next()
$r = $following.Range
$r.Collapse(1)              # wdCollapseStart -- right before the following paragraph's text
$r.InsertBefore("Will research about Prop type and require prop types later `r")
